$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.259.01"
$ws.Range("E2").Value = "  +0.67%  "
$ws.Range("D3").Value = "1.660.03"
$ws.Range("E3").Value = "  +0.57%  "
$ws.Range("E4").Value = "  +0.63%  "
$ws.Range("D5").Value = "'218.45"
$ws.Range("E5").Value = "  +0.19%  "
$ws.Range("D6").Value = "'0.5333"
$ws.Range("E6").Value = "  +0.99%  "
$ws.Range("E8").Value = "  +0.99%  "
$ws.Range("D9").Value = "'0.06348"
$ws.Range("E9").Value = "  +0.64%  "
$ws.Range("D10").Value = "'20.52"
$ws.Range("E10").Value = "  +0.79%  "
$ws.Range("D11").Value = "'0.07831"
$ws.Range("E11").Value = "  +1.18%  "
$ws.Range("D12").Value = "'4.543"
$ws.Range("E12").Value = "  +1.65%  "
$ws.Range("D13").Value = "1.689.22"
$ws.Range("E13").Value = "  +3.01%  "
$ws.Range("D14").Value = "1.888.35"
$ws.Range("E14").Value = "  +0.60%  "
$ws.Range("D15").Value = "'0.5519"
$ws.Range("E15").Value = "  +1.20%  "
$ws.Range("D16").Value = "0.0₅8187"
$ws.Range("E16").Value = "  +0.96%  "
$ws.Range("E17").Value = "  +0.76%  "
$ws.Range("D18").Value = "26.244.49"
$ws.Range("E18").Value = "  +0.55%  "
$ws.Range("E19").Value = "  +0.64%  "
$ws.Range("D20").Value = "'4.652"
$ws.Range("E20").Value = "  +2.39%  "
$ws.Range("D21").Value = "'192.07"
$ws.Range("E21").Value = "  -0.81%  "
$ws.Range("D22").Value = "'10.12"
$ws.Range("E22").Value = "  +0.80%  "
$ws.Range("E23").Value = "  +1.06%  "
$ws.Range("E24").Value = "  +0.59%  "
$ws.Range("D25").Value = "'144.62"
$ws.Range("E25").Value = "  +3.28%  "
$ws.Range("D26").Value = "'0.1229"
$ws.Range("E26").Value = "  -0.73%  "
$ws.Range("D28").Value = "'16.06"
$ws.Range("E28").Value = "  -0.49%  "
$ws.Range("D29").Value = "'1.469"
$ws.Range("E29").Value = "  +2.45%  "
$ws.Range("D30").Value = "'0.05791"
$ws.Range("E30").Value = "  -2.00%  "
$ws.Range("E31").Value = "  -0.02%  "
$ws.Range("D32").Value = "'3.577"
$ws.Range("E32").Value = "  +2.23%  "
$ws.Range("D33").Value = "'3.285"
$ws.Range("E33").Value = "  +1.51%  "
$ws.Range("D34").Value = "'1.609"
$ws.Range("E34").Value = "  +4.13%  "
$ws.Range("D35").Value = "'2.819"
$ws.Range("E35").Value = "  +2.30%  "
$ws.Range("D36").Value = "'0.9571"
$ws.Range("E36").Value = "  +1.48%  "
$ws.Range("E37").Value = "  +0.70%  "
$ws.Range("D38").Value = "'0.5795"
$ws.Range("E38").Value = "  +2.39%  "
$ws.Range("E39").Value = "  -0.02%  "
$ws.Range("D40").Value = "'5.850"
$ws.Range("E40").Value = "  +0.17%  "
$ws.Range("D41").Value = "'0.8537"
$ws.Range("E41").Value = "  +1.02%  "
$ws.Range("D42").Value = "'1.009"
$ws.Range("E42").Value = "  +0.58%  "
$ws.Range("D43").Value = "'104.69"
$ws.Range("E43").Value = "  +3.99%  "
$ws.Range("D44").Value = "1.044.93"
$ws.Range("E44").Value = "  +3.83%  "
$ws.Range("D45").Value = "1.801.41"
$ws.Range("E45").Value = "  +0.44%  "
$ws.Range("E46").Value = "  +0.47%  "
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").Value = "0.0₈105"
$ws.Range("E47").Value = "  -0.78%  "
$ws.Range("B48").Value = "Frax"
$ws.Range("C48").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D48").Value = "'1.012"
$ws.Range("E48").Value = "  +0.52%  "
$ws.Range("E49").Value = "  +1.86%  "
$ws.Range("D50").Value = "'7.950"
$ws.Range("E50").Value = "  +2.03%  "
$ws.Range("D51").Value = "'0.05161"
$ws.Range("E51").Value = "  +0.18%  "
